$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.352.21'
$ws.Range("E2").Value = '  +4.17%  '
$ws.Range("D3").Value = '1.792.07'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5387'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3828'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07559'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.119'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.006'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.156'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.375'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.59%  '
$ws.Range("D16").Value = '1.804.36'
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001069'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06439'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.943'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.48%  '
$ws.Range("D23").Value = '28.394.73'
$ws.Range("E23").Value = '  +4.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.136'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.387'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.64%  '
$ws.Range("D29").Value = '2.010.71'
$ws.Range("E29").Value = '  +2.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.114'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1015'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.725'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.708'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2306'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +14.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06367'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02313'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.145'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.758'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6377'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.155'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.390'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5945'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.674'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.971'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.146'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06896'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.80%  '
